$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row 5, cloned from row 4, so it carries the same
#        cell styles (bold label col, yellow value col, etc.) as rows 2-4.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(5).Insert()

# --- 2. Switch the crystal generator clock from 50 MHz to 49.152 MHz
$ws.Cells.Item(2, 3).Value2 = 49152000

# --- 3. Row 4 becomes "Clock rate ... RX" = 16
$ws.Cells.Item(4, 2).Value2 = "Clock rate при расчёте CIC компенсатора RX"
$ws.Cells.Item(4, 3).Value2 = 16

# --- 4. Row 5 (new) becomes "Clock rate ... TX" = 14
$ws.Cells.Item(5, 2).Value2 = "Clock rate при расчёте CIC компенсатора TX"
$ws.Cells.Item(5, 3).Value2 = 14

# --- 5. Drop the old column E entirely (rows 2-5 had a trailing styled
#        but empty E cell) -- fully clear so the <c> element disappears.
$ws.Cells.Item(2, 5).Clear()
$ws.Cells.Item(3, 5).Clear()
$ws.Cells.Item(4, 5).Clear()
$ws.Cells.Item(5, 5).Clear()

# --- 6. Old row 6 ("Расчёт" / "Практика") shifted down to row 7 by the
#        insert above. Drop the "Практика" side label in column D.
$ws.Cells.Item(7, 4).Clear()

# --- 7. Drop column D all the way down (RX/TX "practice" comparison
#        columns are being removed, only the C column calc stays).
$ws.Cells.Item(8, 4).Clear()
$ws.Cells.Item(9, 4).Clear()

# --- 8. Row 10: column D used to hold the practice 512 value and column
#        E the "powers of two" note; now D10 holds that note text and E
#        is dropped.
$ws.Cells.Item(10, 4).Clear()
$ws.Cells.Item(10, 4).Value2 = "*все коэффициенты должны быть степенью двойки (2,4,8,16,32,…)"
$ws.Cells.Item(10, 5).Clear()

# --- 9. Row 11 label changes, column D dropped.
$ws.Cells.Item(11, 2).Value2 = "Полученная частота дискретизации CIC коспенсатора"
$ws.Cells.Item(11, 4).Clear()

# --- 10. Rows 12-13: just drop column D.
$ws.Cells.Item(12, 4).Clear()
$ws.Cells.Item(13, 4).Clear()

# --- 11. Row 14 becomes the RX CIC-compensator frequency line, column D dropped.
$ws.Cells.Item(14, 2).Value2 = "Частота для CIC компенсатора RX,гц"
$ws.Cells.Item(14, 4).Clear()

# --- 12. New row 15: TX CIC-compensator frequency line.
$ws.Cells.Item(15, 2).Value2 = "Частота для CIC компенсатора RX,гц"
$ws.Cells.Item(15, 3).Formula = "=C8*C5/C10"
$ws.Cells.Item(15, 3).NumberFormat = "0.00"

# --- 13. New row 16: I2S bus frequency line.
$ws.Cells.Item(16, 2).Value2 = "Частота для I2S шины"
$ws.Cells.Item(16, 3).Formula = "=C12*256"

# --- 14. Match the author's final selection.
$ws.Range("C6").Select()
